$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(13)
Write-Output $p.Range.Text
